$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.977.76'
$ws.Range('E2').Value = '  -0.73%  '

# Row 3
$ws.Range('D3').Value = '3.422.90'
$ws.Range('E3').Value = '  -2.64%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').Value = '583.96'
$ws.Range('E5').Value = '  -2.98%  '

# Row 6
$ws.Range('D6').Value = '136.17'
$ws.Range('E6').Value = '  -5.27%  '

# Row 7
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.422.83'
$ws.Range('E7').Value = '  -2.66%  '

# Row 8
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  -3.91%  '

# Row 10
$ws.Range('D10').Value = '0.120'
$ws.Range('E10').Value = '  -9.92%  '

# Row 11
$ws.Range('D11').Value = '7.01'
$ws.Range('E11').Value = '  -10.57%  '

# Row 12
$ws.Range('D12').Value = '0.373'
$ws.Range('E12').Value = '  -7.70%  '

# Row 13
$ws.Range('D13').Value = '4.005.29'
$ws.Range('E13').Value = '  -2.65%  '

# Row 14
$ws.Range('D14').Value = '0.0000177'
$ws.Range('E14').Value = '  -10.75%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.425.13'
$ws.Range('E15').Value = '  -2.48%  '

# Row 16
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '26.20'
$ws.Range('E16').Value = '  -8.08%  '

# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.115'
$ws.Range('E17').Value = '  -1.95%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '64.972.09'
$ws.Range('E18').Value = '  -0.68%  '

# Row 19
$ws.Range('D19').Value = '9.62'
$ws.Range('E19').Value = '  -13.06%  '

# Row 20
$ws.Range('D20').Value = '5.77'
$ws.Range('E20').Value = '  -6.73%  '

# Row 21
$ws.Range('D21').Value = '13.44'
$ws.Range('E21').Value = '  -6.26%  '

# Row 22
$ws.Range('D22').Value = '381.25'
$ws.Range('E22').Value = '  -8.56%  '

# Row 23
$ws.Range('D23').Value = '0.550'
$ws.Range('E23').Value = '  -8.08%  '

# Row 24
$ws.Range('E24').Value = '  -0.08%  '

# Row 25
$ws.Range('D25').Value = '72.16'
$ws.Range('E25').Value = '  -6.97%  '

# Row 26
$ws.Range('D26').Value = '3.553.06'
$ws.Range('E26').Value = '  -2.84%  '

# Row 27
$ws.Range('D27').Value = '0.0000105'
$ws.Range('E27').Value = '  -8.50%  '

# Row 28
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.14%  '

# Row 29
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.09'
$ws.Range('E29').Value = '  -8.87%  '

# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '8.06'
$ws.Range('E30').Value = '  -9.12%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.19'
$ws.Range('E31').Value = '  -10.40%  '

# Row 32
$ws.Range('D32').Value = '3.428.98'
$ws.Range('E32').Value = '  -2.70%  '

# Row 33
$ws.Range('E33').Value = '  -0.03%  '

# Row 34
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.143'
$ws.Range('E34').Value = '  -6.59%  '

# Row 35
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '23.03'
$ws.Range('E35').Value = '  -5.25%  '

# Row 36
$ws.Range('D36').Value = '168.58'
$ws.Range('E36').Value = '  -3.77%  '

# Row 37
$ws.Range('D37').Value = '6.71'
$ws.Range('E37').Value = '  -10.55%  '

# Row 38
$ws.Range('D38').Value = '1.16'
$ws.Range('E38').Value = '  -11.82%  '

# Row 39
$ws.Range('D39').Value = '1.45'
$ws.Range('E39').Value = '  -7.82%  '

# Row 40
$ws.Range('D40').Value = '4.64'
$ws.Range('E40').Value = '  -12.00%  '

# Row 41
$ws.Range('D41').Value = '0.0751'
$ws.Range('E41').Value = '  -7.94%  '

# Row 42
$ws.Range('D42').Value = '0.810'
$ws.Range('E42').Value = '  -5.10%  '

# Row 43
$ws.Range('E43').Value = '  +0.17%  '

# Row 44
$ws.Range('D44').Value = '42.41'
$ws.Range('E44').Value = '  -6.35%  '

# Row 45
$ws.Range('D45').Value = '4.32'
$ws.Range('E45').Value = '  -14.73%  '

# Row 46
$ws.Range('D46').Value = '1.60'
$ws.Range('E46').Value = '  -9.68%  '

# Row 47
$ws.Range('D47').Value = '1.11'
$ws.Range('E47').Value = '  +1.81%  '

# Row 48
$ws.Range('D48').Value = '22.36'
$ws.Range('E48').Value = '  -5.09%  '

# Row 49
$ws.Range('D49').Value = '6.44'
$ws.Range('E49').Value = '  -8.26%  '

# Row 50
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').Value = '2.02'
$ws.Range('E50').Value = '  -14.26%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.162.50'
$ws.Range('E51').Value = '  -7.13%  '
